$d = $word.ActiveDocument

function Replace-AllInRange($range, [string]$findText, [string]$replaceText) {
    $guard = 0
    while ($true) {
        $found = $range.Find.Execute($findText, $true, $false, $false, $false, $false, `
                                      $true, 0, $false, $replaceText, 1)
        if (-not $found) { break }
        $guard = $guard + 1
        if ($guard -gt 200) { break }
    }
}

# --- Table 1 & Table 2: every "DSAI_..." filename mention gets "Kbot_" inserted
# right after the "DSAI_" prefix (rename DSAI_Knowledge_Graph -> DSAI_Kbot_Knowledge_Graph,
# DSAI_Graph.ttl -> DSAI_Kbot_Graph.ttl, etc.)
$t1 = $d.Tables(1)
Replace-AllInRange $t1.Range "DSAI_" "DSAI_Kbot_"

$t2 = $d.Tables(2)
Replace-AllInRange $t2.Range "DSAI_" "DSAI_Kbot_"

# --- Body text after the tables ("Source Code Steps:" section): same rename, but the
# lone "DSAI_Knowledge_Graph_Chatbot" + ".py" (split across two runs) must stay as-is,
# so it is excluded by scoping the body range to stop right before that paragraph.
$bodyStart = $t2.Range.End
$docEnd = $d.Content.End
$bodyRange = $d.Range($bodyStart, $docEnd)
Replace-AllInRange $bodyRange "DSAI_" "DSAI_Kbot_"

# --- Fix up the one occurrence that the original author mistakenly typed with an extra
# space: "...knowledge graph from DSAI_Kbot_KG_Input_Data_Text.csv" should actually read
# "...knowledge graph from DSAI_ Kbot_KG_Input_Data_Text.csv".
Replace-AllInRange $d.Content "from DSAI_Kbot_KG_Input_Data_Text.csv" "from DSAI_ Kbot_KG_Input_Data_Text.csv"
